$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 249.4
$ws.Range("I33").Value = 138.09091
$ws.Range("K33").Value = 138.09091
$ws.Range("M33").Value = 90.90908999999999
$ws.Range("H70").Value = 2012.8334
$ws.Range("I70").Value = 892.3333
$ws.Range("J70").Value = 3133.3333
$ws.Range("K70").Value = 2676.9999
$ws.Range("L70").Value = 9399.999899999999
$ws.Range("M70").Value = -2406.9999
$ws.Range("N70").Value = -9939.999899999999
$ws.Range("H73").Value = 2012.8334
$ws.Range("I73").Value = 892.3333
$ws.Range("J73").Value = 3133.3333
$ws.Range("K73").Value = 2676.9999
$ws.Range("L73").Value = 9399.999899999999
$ws.Range("M73").Value = -1740.9999
$ws.Range("N73").Value = -11271.9999
$ws.Range("H92").Value = 2045.0834
$ws.Range("I92").Value = 2140.0908
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2140.0908
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -892.0907999999999
$ws.Range("N92").Value = -3496
$ws.Range("H98").Value = 1300.04
$ws.Range("I98").Value = 1376.0526
$ws.Range("J98").Value = 1059.3334
$ws.Range("K98").Value = 1376.0526
$ws.Range("L98").Value = 1059.3334
$ws.Range("M98").Value = 121.9474
$ws.Range("N98").Value = -4055.3334
$ws.Range("H113").Value = 3950.182
$ws.Range("I113").Value = 4458.8335
$ws.Range("J113").Value = 3339.8
$ws.Range("K113").Value = 4458.8335
$ws.Range("L113").Value = 3339.8
$ws.Range("M113").Value = -1204.8335
$ws.Range("N113").Value = -9847.799999999999
$ws.Range("H116").Value = 4635
$ws.Range("I116").Value = 5321.4287
$ws.Range("J116").Value = 3834.1667
$ws.Range("K116").Value = 5321.4287
$ws.Range("L116").Value = 3834.1667
$ws.Range("M116").Value = -1879.4287
$ws.Range("N116").Value = -10718.1667
$ws.Range("H117").Value = 24828
$ws.Range("J117").Value = 24828
$ws.Range("L117").Value = 24828
$ws.Range("N117").Value = -34006
$ws.Range("H122").Value = 1300.04
$ws.Range("I122").Value = 1376.0526
$ws.Range("J122").Value = 1059.3334
$ws.Range("K122").Value = 4128.1578
$ws.Range("L122").Value = 3178.0002
$ws.Range("M122").Value = -1678.1578
$ws.Range("N122").Value = -8078.0002
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H135").Value = 3241.0715
$ws.Range("J135").Value = 6630
$ws.Range("L135").Value = 59670
$ws.Range("N135").Value = -64740
$ws.Range("H138").Value = 3005.62
$ws.Range("I138").Value = 1362.1904
$ws.Range("J138").Value = 3442.481
$ws.Range("K138").Value = 4086.5712
$ws.Range("L138").Value = 10327.443
$ws.Range("M138").Value = 1053.4288
$ws.Range("N138").Value = -20607.443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H132").Value = 1794.1316
$ws.Range("I132").Value = 1319.6333
$ws.Range("K132").Value = 3958.8999
$ws.Range("M132").Value = -1428.8999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2955.5
$ws.Range("I10").Value = 2905
$ws.Range("J10").Value = 3006
$ws.Range("K10").Value = 2905
$ws.Range("L10").Value = 3006
$ws.Range("M10").Value = -2765
$ws.Range("N10").Value = -3286
$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 3000
$ws.Range("K75").Value = 3000
$ws.Range("M75").Value = -2064
$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 3000
$ws.Range("K78").Value = 9000
$ws.Range("M78").Value = -4320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1175.5
$ws.Range("I22").Value = 1395.7
$ws.Range("J22").Value = 625
$ws.Range("K22").Value = 1395.7
$ws.Range("L22").Value = 625
$ws.Range("M22").Value = -1045.7
$ws.Range("N22").Value = -1325
$ws.Range("H62").Value = 4650.4
$ws.Range("I62").Value = 6117.5
$ws.Range("J62").Value = 2449.75
$ws.Range("K62").Value = 6117.5
$ws.Range("L62").Value = 2449.75
$ws.Range("M62").Value = -5493.5
$ws.Range("N62").Value = -3697.75
$ws.Range("H65").Value = 4650.4
$ws.Range("I65").Value = 6117.5
$ws.Range("J65").Value = 2449.75
$ws.Range("K65").Value = 30587.5
$ws.Range("L65").Value = 12248.75
$ws.Range("M65").Value = -27467.5
$ws.Range("N65").Value = -18488.75
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52372
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16701246
$ws.Range("I131").Value = 62626492
$ws.Range("J131").Value = 1156.3636
$ws.Range("K131").Value = 187879476
$ws.Range("L131").Value = 3469.0908
$ws.Range("M131").Value = -187874436
$ws.Range("N131").Value = -13549.0908
$ws.Range("H132").Value = 827.1875
$ws.Range("J132").Value = 995
$ws.Range("L132").Value = 8955
$ws.Range("N132").Value = -14015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H97").Value = 1614.5454
$ws.Range("I97").Value = 1576
$ws.Range("K97").Value = 1576
$ws.Range("M97").Value = -1080
$ws.Range("H126").Value = 4636.2
$ws.Range("I126").Value = 4594.222
$ws.Range("J126").Value = 5014
$ws.Range("K126").Value = 13782.666
$ws.Range("L126").Value = 15042
$ws.Range("M126").Value = -11312.666
$ws.Range("N126").Value = -19982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1994.25
$ws.Range("I7").Value = 2237.5
$ws.Range("J7").Value = 1832.0834
$ws.Range("K7").Value = 2237.5
$ws.Range("L7").Value = 1832.0834
$ws.Range("M7").Value = -2125.5
$ws.Range("N7").Value = -2056.0834
$ws.Range("H126").Value = 1994.25
$ws.Range("I126").Value = 2237.5
$ws.Range("J126").Value = 1832.0834
$ws.Range("K126").Value = 6712.5
$ws.Range("L126").Value = 5496.2502
$ws.Range("M126").Value = -4242.5
$ws.Range("N126").Value = -10436.2502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 919
$ws.Range("I126").Value = 631.6667
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 1895.0001
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = 574.9999
$ws.Range("N126").Value = -8990
